$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")
Write-Host $ws.Name
Write-Host $ws.Range("A97").Value
